# Gages.xlsx - "Marys" sheet update
# Insert a new gage row (Outlet of Muddy Creek into the Marys) above the
# existing "Marys River outlet into the Willamette" row, and tweak the
# number format / selection on the top data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marys")
$ws.Activate()

# Insert a new row at row 3 (pushes the existing rows 3-9 down to 4-10)
$ws.Rows.Item(3).Insert()

# Populate the new row with the Muddy Creek outlet gage data
$ws.Range("B3").Value = "Outlet of Muddy Creek into the Marys"
$ws.Range("D3").Value = 23762959
$ws.Range("E3").Value = 31856.486400000002

# Row 2's computed-area cell (M2) switches from a whole-number format to a
# 2-decimal number format
$ws.Range("M2").NumberFormat = "0.00"

# Update the visible selection to M2
[void]$ws.Range("M2").Select()
